$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly documented "Descripción" values in column B for the
# rows that previously had no description (order chosen to match the
# shared-string table layout produced by the original author's edit).
$ws.Range("B7").Value  = "Fecha en la que se facturó"
$ws.Range("B27").Value = "Fecha y hora de la transacción (repetido)"
$ws.Range("B28").Value = "Fecha de cuando la empresa empieza a trabajar el pedido"
$ws.Range("B31").Value = "Id interno de la dirección de entrega"
$ws.Range("B23").Value = "Código interno del proveedor"
$ws.Range("B6").Value  = "Número de lote"
$ws.Range("B26").Value = "Descripción del ecommerce"
$ws.Range("B20").Value = "Sucursal de donde salió el producto"
$ws.Range("B21").Value = "Sucursal a donde se envió el producto por última vez"
$ws.Range("B29").Value = "Fecha en la que se facturó"

# The header row no longer needs the bold style.
$ws.Range("A1:B1").Font.Bold = $false

# Move the view back to the top of the sheet, matching where the author
# left the cursor after finishing the dictionary.
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 10
